{"js": "// Replace each two-digit multiplication expression with its updated\n// version, matching the exact old text so only the intended table cell\n// is touched (every expression in the document is unique).\nconst replacements = [\n  [\"83\u00d731=2573\", \"77\u00d728=2156\"],\n  [\"17\u00d759=1003\", \"72\u00d730=2160\"],\n  [\"27\u00d752=1404\", \"74\u00d787=6438\"],\n  [\"47\u00d737=1739\", \"79\u00d726=2054\"],\n  [\"80\u00d778=6240\", \"48\u00d741=1968\"],\n  [\"83\u00d717=1411\", \"74\u00d784=6216\"],\n  [\"56\u00d781=4536\", \"83\u00d712=996\"],\n  [\"59\u00d721=1239\", \"14\u00d787=1218\"],\n  [\"25\u00d725=625\", \"71\u00d799=7029\"],\n  [\"56\u00d750=2800\", \"67\u00d732=2144\"],\n  [\"69\u00d797=6693\", \"50\u00d728=1400\"],\n  [\"80\u00d795=7600\", \"31\u00d739=1209\"],\n  [\"65\u00d754=3510\", \"70\u00d777=5390\"],\n  [\"80\u00d780=6400\", \"97\u00d713=1261\"],\n  [\"50\u00d780=4000\", \"49\u00d712=588\"],\n  [\"19\u00d782=1558\", \"98\u00d729=2842\"],\n  [\"51\u00d723=1173\", \"74\u00d711=814\"],\n  [\"65\u00d718=1170\", \"13\u00d793=1209\"],\n  [\"20\u00d797=1940\", \"30\u00d716=480\"],\n  [\"32\u00d738=1216\", \"20\u00d788=1760\"],\n  [\"44\u00d780=3520\", \"32\u00d789=2848\"],\n  [\"73\u00d745=3285\", \"11\u00d713=143\"],\n  [\"72\u00d766=4752\", \"16\u00d719=304\"],\n  [\"86\u00d799=8514\", \"58\u00d725=1450\"],\n  [\"33\u00d770=2310\", \"11\u00d778=858\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication expression with its updated\n# version. Every expression in the document is unique, so a plain\n# Find/Replace (MatchWholeWord off, MatchCase on) targets exactly the\n# intended table cell each time.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"83\u00d731=2573\"; New = \"77\u00d728=2156\" }\n    @{ Old = \"17\u00d759=1003\"; New = \"72\u00d730=2160\" }\n    @{ Old = \"27\u00d752=1404\"; New = \"74\u00d787=6438\" }\n    @{ Old = \"47\u00d737=1739\"; New = \"79\u00d726=2054\" }\n    @{ Old = \"80\u00d778=6240\"; New = \"48\u00d741=1968\" }\n    @{ Old = \"83\u00d717=1411\"; New = \"74\u00d784=6216\" }\n    @{ Old = \"56\u00d781=4536\"; New = \"83\u00d712=996\" }\n    @{ Old = \"59\u00d721=1239\"; New = \"14\u00d787=1218\" }\n    @{ Old = \"25\u00d725=625\"; New = \"71\u00d799=7029\" }\n    @{ Old = \"56\u00d750=2800\"; New = \"67\u00d732=2144\" }\n    @{ Old = \"69\u00d797=6693\"; New = \"50\u00d728=1400\" }\n    @{ Old = \"80\u00d795=7600\"; New = \"31\u00d739=1209\" }\n    @{ Old = \"65\u00d754=3510\"; New = \"70\u00d777=5390\" }\n    @{ Old = \"80\u00d780=6400\"; New = \"97\u00d713=1261\" }\n    @{ Old = \"50\u00d780=4000\"; New = \"49\u00d712=588\" }\n    @{ Old = \"19\u00d782=1558\"; New = \"98\u00d729=2842\" }\n    @{ Old = \"51\u00d723=1173\"; New = \"74\u00d711=814\" }\n    @{ Old = \"65\u00d718=1170\"; New = \"13\u00d793=1209\" }\n    @{ Old = \"20\u00d797=1940\"; New = \"30\u00d716=480\" }\n    @{ Old = \"32\u00d738=1216\"; New = \"20\u00d788=1760\" }\n    @{ Old = \"44\u00d780=3520\"; New = \"32\u00d789=2848\" }\n    @{ Old = \"73\u00d745=3285\"; New = \"11\u00d713=143\" }\n    @{ Old = \"72\u00d766=4752\"; New = \"16\u00d719=304\" }\n    @{ Old = \"86\u00d799=8514\"; New = \"58\u00d725=1450\" }\n    @{ Old = \"33\u00d770=2310\"; New = \"11\u00d778=858\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute(\n        $pair.Old,   # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $pair.New,   # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        throw \"Could not find text: $($pair.Old)\"\n    }\n}\n"}
